$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename sheet "My Series" -> "Data"
$ws.Name = "Data"

# 2. Cell text edits
$ws.Range("A11").Value = "Function Information"
$ws.Range("B20").Value = -0.3785388757796309

# 3. Number format change for the decimal stats block (numFmtId 166: 0.000 -> ###0.000)
$ws.Range("B27:B36").NumberFormat = "###0.000"

# 4. Update the hidden CEIC add-in metadata comment on A1 to the refreshed payload
$comment = $ws.Range("A1").Comment
$comment.Text("8RwAAB+LCAAAAAAAAAOlWVtvG8cV/isLPrUP5C4vjkVhtAEvkkqUFAWRqqy8FMvdkTj1cofdmZXEtxRJESOtgQJNXLRNWrRA0aJAkwBNgyZu4f8SmLL71L/QM5e9kZTNdQ3B2jm3OXPmzDnfjNDbNzPfuMIhIzTYK1UrVsnAgUs9ElzulSJ+Ua6+VXrbRvs3LvaPndCZYQ7CBmgFbPeGkb3SlPP5rmleX19XrusVGl6aNcuqmg8G/ZE7xTOnTALGncDFpUTLe71WyUYdbzbA3PEc7ijNvVJv1Kt0MHG7QBs4gXOJw0o7YiTAjO0HnHCCmdAMscNxpzv4gVqYXau8Vakic42eSrYj4ntKLiep6FoOpsVjMsN2zarulK1m2aqPreZuo7Fbq1V2mtV3YsVEEPUdxkc4vCKuJIy4M5tLdatp1a1mo1GrIXOjENhKA2Cjoe+d4CvCsNfBvs8KRcTUG9hyOay6WDAtZGZ0taE3d+EwdObTMeE+LubGyWnbmATal9SIjQ5oiF2I3xu5dISvh6EO63jeB+54SkK+6DqLwrZOGQ6HcxGkYqo26tKAt3wc8tM57DX2IBWAYfMwwsi8g5kqdQlz4ZsEEfbsC8dnWaUcE53R8CGbOy4+gnNsChvXgU8dDxKOE8aJm066xkDHIZ2DRZi8TX3vAKxq4Q2MxHIvgBCLaduUPky928REclfl/sKezhwei6/R0WhKr4eBvxhFE+aGZIK9bjuW3shD4kBq7U7EOJ2BFykJKVqGMhiYC/gHZ3CVg7rYJTPHP/Yhjsyug6EcAbUiTi8I71A/mgUsdmuFis5gUWN8kywyGaMh7G8g4k6DXhDLq0hvZOUVTuh1Muc6Q8YhQ24xN97xdcaqcBdo8Q6uc+SmiFUeEB9aRHY7MtR8YoymGPONWaE4SFTDA9F07PbiKJpN4IRN4JhdyVkZMlM+glSFdAe/bAsaSVn+jC1rV/6AHwkb7Qfe3XIxE8F0mbnsKvBWSAjW5Ld9J3gI1DPCp0eteC0bOEhF4E75dR6Cwzv3nYUkJ1HK0lAvcP3Iw6om9IILmaLCN7Wpd7LRGqkPx9xGTrAYL+ZQmhnZ5fCxV4Jmvct4CHCgZLs0Cni4EMUDmVr0dTosmgRyAsffWucixD+OAIUsDqLA7VBv+9k8FZ3TgPDtPaRRqCri9ioyeqI4RqyLRZmRdX9rfbfImlhYSHwW4BkNiLt9tCHIwnvvDRbC4lO1tQZW52treR9au+p84qxvrRYChIReV2iaFmPUJTJZ9fHwMvrmHUemiy+cyAf4xqHLXia1d5WMWuzhqkyWhE5DP66AtgDHDNCx680qLuAHgQArLp0Jggmg9GyEzKy8AEEu3g8u+05wGQHMSOrKKj2pv6JFjkMnYGI5CapYKcWbhVBcpxTasVXxGkYyEVTxosBF5oocGuPZnIaOP4DAkAOddhoyARoZOHyqR9DbfOzGQTZT1UQr71ns+OvEZJNSyxAHXpfJFaIUEmtRODyVSWlIrHIAx9LvOD6ZhKqqxq18Ew82LMWHcf0ViyuIFeM9gLsYdN/v44UA6OlA02XKVmOGSmBRSO3RSWOn1mjugJIcI7niDvhEfSLxpXFI4ToYzEDPAGSMATnuGgsOR0EFRzegAnpZDXQQ13Y5/QBSNEfI8wFiXBLoHGtyCSOVtwcAPKf+IiOpVtenLgguf/Fo+fVXtx89W37zye1HXy1/+qf//us3y6+/XD769OUHf3/+9EO1PiWMxs7Ex9KhcXtnx6o3ILUSEhLxNCUg9iKXS9r5ucTByRjp65wcdPZ7ncN+W5aQhBirqy5iipvigkbpcKQWISeSu2jGm69E7HFckvQ4x810JVvc3K5wXjrLv0tRxeLF0z+/ePq3O7V1wFJ4VW0269ugL7gDN8pW7Q701c+VfS18r1yrZYRXZNCJqvlJnHqeXa9azerOfaualG8vyeBNQqssbWnsXJoreorUUYgoSYHsOGbKvB9jxhO2OgmZgQri7T9+9vLzj3NSOrqakrcCzknkIiYz44E0fXQyNkbD05POvjHeH4k8SXkZOWX8FcJ69uQ8ZZNm/bTlUmoTVZl7peLqjIchjeZrlSKlbpDcWDPWuSuaMjbrrqW8DfJ61z55tElBr6SbAtLkCSNLQzmOImX4unY9/snzb549f/b+i3//6vbJ58vHf1z+872cHT1bcr+H1IUDkh0mmQxVTHeNFQo6G8mYPrR+mOkSmiiuQ8eUBJzZtXsNeRXSQwS6VWFO/ka9GXQuaVmGDegrFPQ9h+3fcH1Y7SNk5gng6NyBpknTK2RCUHU5De9/Pv3d7W+/vH3yxcsP/rr88C/Lnz958fT3Lz/7gzpJtx9/cfv4M125V4u79EVcTBWWM+RLh2uIE2aIFmx8++4vjYByA5CDEckq8+27v84YE45KjJFaBmSWOJJ3YU00qyz0jIwriQ85vURF9fGOaEv1REI3JjonbjrJO2VhCuClIRnf6Y3LEcMGBVD0XVhJXjhV3lZPq6g2eXzfqlVrmqu8EUuYOCwT+kOfThzfiBnyHWFFJKf1aoVUVs532B+2W/1URDkxDD0cijRUHyhGhqJN9Fg8ilMtQwEu4Dc38sXTz5rYOiuxnClnpn5DuWh54kF984tDTgJ1ojBUKCfQr/KjaA6YNn5su5svHyAzMPZIQc4ssE3HvW6eD+MMF5pbni0Iki8rlGapatVj4rlGodIjEZp0CLzcoyWEQz+8K/R0RWBLTFF49sOQhhurT8qJxQYAiKGimGnEExm5pwo8e+lexYS44iUf6gKnV0i72Me82Ku0mWoPAP2+qS7sfVHVHhv6ng5msRtEEpbUQPZpXiTK//syr5KtFYYAlsQ7XuGn9Pj+eQLX1oLeqKVIRXGRg9n1W/cBCRl/ICqB/lKU84RyrlDnA3FvUh9yfG7X7ykCCJhZ62bOzfjocvUHEur3yYwUvN1Z8fnOG4FYzucKlvWKZYpoLUf4BkBjxgIUxcmPoG2o55Ai1lTCQi1N9MUTJCOXU17UsfsTB3t4YpXdCa6VG561U25iXC9Xq/C/49ZqlnVPPGBq41A5CL4uOIkZb1j6h0v7fx++edDxHAAA")
